$d = $word.ActiveDocument

# 1. Fix the GD term end date: "2018 - Present" -> "2018 - April 2019"
$d.Content.Find.Execute(" 2018 – Present", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " 2018 – April 2019", 2)
